$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 6
$ws.Range("H6").Value = 669.8
$ws.Range("I6").Value = 669.8
$ws.Range("K6").Value = 2009.4
$ws.Range("M6").Value = -1897.4

# ALC row 24
$ws.Range("H24").Value = 1618
$ws.Range("I24").Value = 1618
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 4854
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = -4684
$ws.Range("N24").ClearContents()

# ALC row 112
$ws.Range("H112").Value = 1328.5227
$ws.Range("J112").Value = 1340.814
$ws.Range("L112").Value = 4022.442
$ws.Range("N112").Value = -6238.442

# ALC row 116
$ws.Range("H116").Value = 5902.5
$ws.Range("I116").Value = 1918.5333
$ws.Range("K116").Value = 1918.5333
$ws.Range("M116").Value = 1523.4667

# ALC row 138
$ws.Range("H138").Value = 3537.7058
$ws.Range("I138").Value = 2830.125
$ws.Range("J138").Value = 4166.6665
$ws.Range("K138").Value = 8490.375
$ws.Range("L138").Value = 12499.9995
$ws.Range("M138").Value = -3350.375
$ws.Range("N138").Value = -22779.9995

$ws = $wb.Worksheets.Item("ARM")
# ARM row 24
$ws.Range("H24").Value = 33333.332
$ws.Range("J24").Value = 33333.332
$ws.Range("L24").Value = 33333.332
$ws.Range("N24").Value = -34081.332

# ARM row 32
$ws.Range("H32").Value = 3773.366
$ws.Range("I32").Value = 3653.9565
$ws.Range("K32").Value = 3653.9565
$ws.Range("M32").Value = -3366.9565

# ARM row 48
$ws.Range("H48").Value = 74800
$ws.Range("J48").Value = 74800
$ws.Range("L48").Value = 74800
$ws.Range("N48").Value = -75568

# ARM row 61
$ws.Range("H61").Value = 1945.909
$ws.Range("I61").Value = 2010.5
$ws.Range("K61").Value = 2010.5
$ws.Range("M61").Value = -1798.5

# ARM row 100
$ws.Range("H100").Value = 33333.332
$ws.Range("J100").Value = 33333.332
$ws.Range("L100").Value = 33333.332
$ws.Range("N100").Value = -35497.332

# ARM row 136
$ws.Range("H136").Value = 1945.909
$ws.Range("I136").Value = 2010.5
$ws.Range("K136").Value = 6031.5
$ws.Range("M136").Value = -3481.5

# ARM row 137
$ws.Range("H137").Value = 38672.57
$ws.Range("J137").Value = 38672.57
$ws.Range("L137").Value = 38672.57
$ws.Range("N137").Value = -48872.57

$ws = $wb.Worksheets.Item("BSM")
# BSM row 99
$ws.Range("H99").Value = 4191.4287
$ws.Range("I99").Value = 1986.6666
$ws.Range("K99").Value = 1986.6666
$ws.Range("M99").Value = -488.6666

# BSM row 137
$ws.Range("H137").Value = 55715
$ws.Range("J137").Value = 55715
$ws.Range("L137").Value = 55715
$ws.Range("N137").Value = -65915

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 231428.08
$ws.Range("I31").Value = 1040631.7
$ws.Range("J31").Value = 2740.1086
$ws.Range("K31").Value = 1040631.7
$ws.Range("L31").Value = 2740.1086
$ws.Range("M31").Value = -1040336.7
$ws.Range("N31").Value = -3330.1086

# CRP row 34
$ws.Range("H34").Value = 231428.08
$ws.Range("I34").Value = 1040631.7
$ws.Range("J34").Value = 2740.1086
$ws.Range("K34").Value = 1040631.7
$ws.Range("L34").Value = 2740.1086
$ws.Range("M34").Value = -1040429.7
$ws.Range("N34").Value = -3144.1086

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 495273.88
$ws.Range("I5").Value = 398.5
$ws.Range("J5").Value = 703642.5
$ws.Range("K5").Value = 1195.5
$ws.Range("L5").Value = 2110927.5
$ws.Range("M5").Value = -1083.5
$ws.Range("N5").Value = -2111151.5

# CUL row 68
$ws.Range("H68").Value = 4995.3228
$ws.Range("I68").Value = 1757.4286
$ws.Range("J68").Value = 5939.7085
$ws.Range("K68").Value = 5272.2858
$ws.Range("L68").Value = 17819.1255
$ws.Range("M68").Value = -4461.2858
$ws.Range("N68").Value = -19441.1255

# CUL row 71
$ws.Range("H71").Value = 4995.3228
$ws.Range("I71").Value = 1757.4286
$ws.Range("J71").Value = 5939.7085
$ws.Range("K71").Value = 15816.8574
$ws.Range("L71").Value = 53457.3765
$ws.Range("M71").Value = -11760.8574
$ws.Range("N71").Value = -61569.3765

# CUL row 109
$ws.Range("H109").Value = 2774.111
$ws.Range("I109").Value = 593.4
$ws.Range("J109").Value = 5500
$ws.Range("K109").Value = 1780.2
$ws.Range("L109").Value = 16500
$ws.Range("M109").Value = -740.1999999999998
$ws.Range("N109").Value = -18580

# CUL row 112
$ws.Range("H112").Value = 585249.7
$ws.Range("J112").Value = 1170000
$ws.Range("L112").Value = 3510000
$ws.Range("N112").Value = -3512216

# CUL row 113
$ws.Range("H113").Value = 3125739.8
$ws.Range("I113").Value = 636.7241
$ws.Range("J113").Value = 11364648
$ws.Range("K113").Value = 1910.1723
$ws.Range("L113").Value = 34093944
$ws.Range("M113").Value = 259.8276999999998
$ws.Range("N113").Value = -34098284

# CUL row 115
$ws.Range("H115").Value = 4366.6665
$ws.Range("I115").Value = 550
$ws.Range("J115").Value = 12000
$ws.Range("K115").Value = 1650
$ws.Range("L115").Value = 36000
$ws.Range("M115").Value = -475
$ws.Range("N115").Value = -38350

# CUL row 118
$ws.Range("H118").Value = 2768.5557
$ws.Range("I118").Value = 702.4286
$ws.Range("J118").Value = 10000
$ws.Range("K118").Value = 2107.2858
$ws.Range("L118").Value = 30000
$ws.Range("M118").Value = -864.2857999999997
$ws.Range("N118").Value = -32486

# CUL row 122
$ws.Range("H122").Value = 2297.982
$ws.Range("I122").Value = 743.6667
$ws.Range("J122").Value = 3054.1353
$ws.Range("K122").Value = 6693.0003
$ws.Range("L122").Value = 27487.2177
$ws.Range("M122").Value = -4243.0003
$ws.Range("N122").Value = -32387.2177

# CUL row 129
$ws.Range("H129").Value = 1792.6316
$ws.Range("I129").Value = 1321.6666
$ws.Range("J129").Value = 2600
$ws.Range("K129").Value = 3964.9998
$ws.Range("L129").Value = 7800
$ws.Range("M129").Value = 1035.0002
$ws.Range("N129").Value = -17800

# CUL row 132
$ws.Range("H132").Value = 2843.5
$ws.Range("I132").Value = 1103.6666
$ws.Range("J132").Value = 4583.3335
$ws.Range("K132").Value = 9932.999400000001
$ws.Range("L132").Value = 41250.0015
$ws.Range("M132").Value = -7402.999400000001
$ws.Range("N132").Value = -46310.0015

# CUL row 133
$ws.Range("H133").Value = 8587.25
$ws.Range("I133").Value = 12374.5
$ws.Range("K133").Value = 37123.5
$ws.Range("M133").Value = -32063.5

# CUL row 134
$ws.Range("H134").Value = 3414.24
$ws.Range("I134").Value = 2065.8462
$ws.Range("J134").Value = 4875
$ws.Range("K134").Value = 6197.5386
$ws.Range("L134").Value = 14625
$ws.Range("M134").Value = -1127.5386
$ws.Range("N134").Value = -24765

# CUL row 135
$ws.Range("H135").Value = 495273.88
$ws.Range("I135").Value = 398.5
$ws.Range("J135").Value = 703642.5
$ws.Range("K135").Value = 3586.5
$ws.Range("L135").Value = 6332782.5
$ws.Range("M135").Value = -1051.5
$ws.Range("N135").Value = -6337852.5

# CUL row 136
$ws.Range("H136").Value = 3472.2354
$ws.Range("I136").Value = 1346.8572
$ws.Range("J136").Value = 4960
$ws.Range("K136").Value = 4040.5716
$ws.Range("L136").Value = 14880
$ws.Range("M136").Value = 1059.4284
$ws.Range("N136").Value = -25080

# CUL row 139
$ws.Range("H139").Value = 1060.9524

$ws = $wb.Worksheets.Item("ARM")
# ARM row 4
$ws.Range("H4").Value = 17250
$ws.Range("J4").Value = 17250
$ws.Range("L4").Value = 17250
$ws.Range("N4").Value = -17474

$ws = $wb.Worksheets.Item("GSM")
# GSM row 6
$ws.Range("H6").Value = 11904
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 11904
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 11904
$ws.Range("N6").Value = -12130
$ws.Range("M6").ClearContents()

# GSM row 16
$ws.Range("H16").Value = 11904
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 11904
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 11904
$ws.Range("N16").Value = -12404
$ws.Range("M16").ClearContents()

# GSM row 46
$ws.Range("H46").Value = 24093.572
$ws.Range("J46").Value = 24004.615
$ws.Range("L46").Value = 24004.615
$ws.Range("N46").Value = -24316.615

# GSM row 107
$ws.Range("H107").Value = 8547940
$ws.Range("I107").Value = 342.85715
$ws.Range("J107").Value = 18520136
$ws.Range("K107").Value = 342.85715
$ws.Range("L107").Value = 18520136
$ws.Range("M107").Value = 1577.14285
$ws.Range("N107").Value = -18523976

# GSM row 113
$ws.Range("H113").Value = 1913.125
$ws.Range("I113").Value = 2158.4
$ws.Range("J113").Value = 1504.3334
$ws.Range("K113").Value = 2158.4
$ws.Range("L113").Value = 1504.3334
$ws.Range("M113").Value = 11.59999999999991
$ws.Range("N113").Value = -5844.3334

# GSM row 122
$ws.Range("H122").Value = 6066.9
$ws.Range("I122").Value = 5124.875
$ws.Range("K122").Value = 15374.625
$ws.Range("M122").Value = -12924.625

# GSM row 137
$ws.Range("H137").Value = 30308
$ws.Range("J137").Value = 45770
$ws.Range("L137").Value = 45770
$ws.Range("N137").Value = -55970

$ws = $wb.Worksheets.Item("WVR")
# WVR row 100
$ws.Range("H100").Value = 519.6
$ws.Range("J100").Value = 600
$ws.Range("L100").Value = 1200
$ws.Range("N100").Value = -2282
